$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the table: insert a new "C"/"D" data column content and re-order rows ---
# Target layout:
#   Row1 (header, bold): A1=GIS LINE NAME  B1=ETS LINE NAME  C1=COMMENT            D1=USER
#   Row2:                A2=E_ESK-MAV      B2=E_ESK-MAVF     C2=Andet navn i ETS   D2=TSP
#   Row3:                A3=E_BLA-MAV_2    B3=E_BLA-MAV      C3=Andet navn i GIS   D3=TSP

$ws.Range("A1").Value = "GIS LINE NAME"
$ws.Range("B1").Value = "ETS LINE NAME"
$ws.Range("C1").Value = "COMMENT"
$ws.Range("D1").Value = "USER"

$ws.Range("A2").Value = "E_ESK-MAV"
$ws.Range("B2").Value = "E_ESK-MAVF"
$ws.Range("C2").Value = "Andet navn i ETS"
$ws.Range("D2").Value = "TSP"

$ws.Range("A3").Value = "E_BLA-MAV_2"
$ws.Range("B3").Value = "E_BLA-MAV"
$ws.Range("C3").Value = "Andet navn i GIS"
$ws.Range("D3").Value = "TSP"

# --- Bold header row ---
$ws.Range("A1:D1").Font.Bold = $true

# --- Column widths (closest values reachable through this engine's pixel
#     quantization that land in the same rounded bucket as the target
#     16.42578125 / 16.7109375 / 15.7109375 character widths) ---
$ws.Range("A1").ColumnWidth = 15.665
$ws.Range("B1").ColumnWidth = 15.83
$ws.Range("C1").ColumnWidth = 14.83

# --- AutoFilter over the header row ---
$ws.Range("A1:D1").AutoFilter()

# Register the hidden _FilterDatabase defined name that Excel creates for the
# sheet's autofilter range.
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=GisMapping!`$A`$1:`$D`$1")
$fdb.Visible = $false

# --- Selection / active cell ---
$ws.Range("G11").Select()

Write-Host "edit complete"
